# Run 5 for span 10: extend the Q-factor sweep table (rows 73-117)
# with the newly completed simulation results (A = detuning, B = run,
# C = Q (dB); D-H are constant per-run metadata).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
  @{Row=73; A=2; B=5; C=-1.4611},
  @{Row=74; A=-9; B=6; C=4.6033},
  @{Row=75; A=-8; B=6; C=5.2784},
  @{Row=76; A=-7; B=6; C=5.8279},
  @{Row=77; A=-6; B=6; C=6.3401},
  @{Row=78; A=-5; B=6; C=6.6346},
  @{Row=79; A=-4; B=6; C=6.7814},
  @{Row=80; A=-3; B=6; C=6.4291},
  @{Row=81; A=-2; B=6; C=5.3739},
  @{Row=82; A=-1; B=6; C=5.1433},
  @{Row=83; A=0; B=6; C=3.2795},
  @{Row=84; A=1; B=6; C=0.64904},
  @{Row=85; A=2; B=6; C=-2.2381},
  @{Row=86; A=-9; B=7; C=4.5494},
  @{Row=87; A=-8; B=7; C=5.1638},
  @{Row=88; A=-7; B=7; C=5.7085},
  @{Row=89; A=-6; B=7; C=6.2424},
  @{Row=90; A=-5; B=7; C=6.5467},
  @{Row=91; A=-4; B=7; C=6.7052},
  @{Row=92; A=-3; B=7; C=6.5723},
  @{Row=93; A=-2; B=7; C=5.7703},
  @{Row=94; A=-1; B=7; C=5.1051},
  @{Row=95; A=0; B=7; C=3.7594},
  @{Row=96; A=1; B=7; C=2.2049},
  @{Row=97; A=2; B=7; C=-0.62087},
  @{Row=98; A=-9; B=8; C=4.5279},
  @{Row=99; A=-8; B=8; C=5.0744},
  @{Row=100; A=-7; B=8; C=5.5962},
  @{Row=101; A=-6; B=8; C=6.1046},
  @{Row=102; A=-5; B=8; C=6.5001},
  @{Row=103; A=-3; B=8; C=6.7209},
  @{Row=104; A=-2; B=8; C=6.3868},
  @{Row=105; A=1; B=8; C=1.7325},
  @{Row=106; A=-9; B=9; C=4.44},
  @{Row=107; A=-5; B=9; C=6.4409},
  @{Row=108; A=-3; B=9; C=6.7729},
  @{Row=109; A=-2; B=9; C=6.5031},
  @{Row=110; A=-1; B=9; C=5.6733},
  @{Row=111; A=-8; B=10; C=4.9098},
  @{Row=112; A=-7; B=10; C=5.4203},
  @{Row=113; A=-4; B=10; C=6.6888},
  @{Row=114; A=-3; B=10; C=6.8123},
  @{Row=115; A=0; B=10; C=4.5922},
  @{Row=116; A=1; B=10; C=4.6839},
  @{Row=117; A=2; B=10; C=3.0886}
)

foreach ($item in $newData) {
  $r = $item.Row
  $ws.Cells.Item($r, 1).Value = $item.A
  $ws.Cells.Item($r, 2).Value = $item.B
  $ws.Cells.Item($r, 3).Value = $item.C
  $ws.Cells.Item($r, 4).Value = 0
  $ws.Cells.Item($r, 5).Value = 250
  $ws.Cells.Item($r, 6).Value = 112
  $ws.Cells.Item($r, 7).Value = 0.16
  $ws.Cells.Item($r, 8).Value = 0.158
}
